$wb = $excel.ActiveWorkbook

# Overview sheet: bump the "Latest HO Xliff Generate Date" for the
# 43a2342d-... row (row 2) now that a fresh handback report was generated.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 02:45:04"

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback
# DateTime for the 43a2342d-... row (row 2).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 02:44:56"
$wsZhCn.Range("K2").Value = "2016-08-19 02:45:28"

# de-de sheet: update Correspond Handoff Datetime / Correspond Handback
# DateTime for the 43a2342d-... row (row 2).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 02:45:04"
$wsDeDe.Range("K2").Value = "2016-08-19 02:45:35"
